# The "SNS VCF" shape on slide 1 originally holds the trailing "VCF " text
# as a single run. PowerPoint re-split that run into "VCF" + " " (the
# space survives as its own run) once the slide was touched again while
# starting to update the deck for Mutect VCFs. Reproduce the split here.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(7)   # "Freeform 51" - contains the "SNS VCF" label
$tr  = $shp.TextFrame.TextRange

# Sanity check anchor text before editing (helps debugging if shape index
# or layout ever drifts): "SNS VCF " (8 characters).
# Characters(Start, Length) is 1-based, matching the VBA TextRange API.

$vcfRun   = $tr.Characters(5, 3)   # "VCF"
$spaceRun = $tr.Characters(8, 1)   # trailing " "

# Re-assert the text on each half so the single run backing "VCF " is
# split into two runs: "VCF" and " ".
$vcfRun.Text   = "VCF"
$spaceRun.Text = " "
